$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")
$ws.Range("A3").Value = "mngr523220"
$ws.Range("B3").Value = "gynUnYd"
$ws.Range("A4:B4").Value = ""
$ws.Range("A3").Select()
